# Insert 9 new historical rows (2019-11-18 .. 2019-11-28) before the
# existing row 890 (2019-11-29), pushing all subsequent rows down by 9.
# This mirrors the commit "Added history starting from 2009-March-31
# historical data to all stocks" for this sheet's data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 890..962 down to 899..971 by inserting 9 blank rows at 890.
$ws.Range("A890:A898").EntireRow.Insert()

# New data rows, in column order: timestamp, date, id, name, open, high, low, close, vol
$newRows = @(
    @(1574035200, "2019-11-18", "5275", "MYNEWS", 1.27, 1.28, 1.27, 1.28, 68000),
    @(1574121600, "2019-11-19", "5275", "MYNEWS", 1.27, 1.29, 1.26, 1.27, 123000),
    @(1574208000, "2019-11-20", "5275", "MYNEWS", 1.27, 1.27, 1.27, 1.27, 88600),
    @(1574294400, "2019-11-21", "5275", "MYNEWS", 1.27, 1.27, 1.26, 1.26, 104500),
    @(1574380800, "2019-11-22", "5275", "MYNEWS", 1.26, 1.27, 1.26, 1.26, 163200),
    @(1574640000, "2019-11-25", "5275", "MYNEWS", 1.28, 1.28, 1.26, 1.26, 152000),
    @(1574726400, "2019-11-26", "5275", "MYNEWS", 1.26, 1.26, 1.25, 1.26, 170600),
    @(1574812800, "2019-11-27", "5275", "MYNEWS", 1.26, 1.26, 1.25, 1.25, 32300),
    @(1574899200, "2019-11-28", "5275", "MYNEWS", 1.25, 1.25, 1.23, 1.24, 594200)
)

$startRow = 890
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $ws.Cells.Item($r, 6).Value2 = $row[5]
    $ws.Cells.Item($r, 7).Value2 = $row[6]
    $ws.Cells.Item($r, 8).Value2 = $row[7]
    $ws.Cells.Item($r, 9).Value2 = $row[8]
}
